# Update "want to go" counts (F2, F4) on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 328
    $ws.Range("F4").Value = 61
}
